$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("noZ")
$ws.Activate()

# Update selection from G11 to E11
$ws.Range("E11").Select()

# Update formula in O2 (unique formula, not shared)
$ws.Range("O2").Formula = "=10^4*N2*(F2/(H2/1000))"

# Update the rest of column O (O3:O9) with the new formula in one shot so
# Excel fills in the relative references for each row.
$ws.Range("O3:O9").Formula = "=10^4*N3*(F3/(H3/1000))"
